$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-7) are being reshuffled: each row's data (date, quality,
# volume, prices, origin, price/kg) is replaced with the data that
# previously belonged to a different row, per the mapping below.
# (columns A,B,C,E,F,G,H,I,J,K,Q,T are identical across rows and unaffected)

$newRow2 = @{ D = 44252; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; R = "Región Metropolitana";    S = 750 }
$newRow3 = @{ D = 44257; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";    S = 806 }
$newRow4 = @{ D = 44250; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";    S = 806 }
$newRow5 = @{ D = 45072; L = "Segunda"; M = 100; N = 16000; O = 16000; P = 16000; R = "Provincia de Chacabuco"; S = 889 }
$newRow6 = @{ D = 45072; L = "Segunda"; M = 100; N = 17000; O = 17000; P = 17000; R = "Provincia de Limarí";    S = 944 }
$newRow7 = @{ D = 44253; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";    S = 806 }

$rows = @{
    2 = $newRow2
    3 = $newRow3
    4 = $newRow4
    5 = $newRow5
    6 = $newRow6
    7 = $newRow7
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("D$r").Value = $data.D
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("O$r").Value = $data.O
    $ws.Range("P$r").Value = $data.P
    $ws.Range("R$r").Value = $data.R
    $ws.Range("S$r").Value = $data.S
}
